$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.003.53"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").Value = "1.920.44"
$ws.Range("E3").Value = "  +0.51%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'325.18"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "'0.4587"
$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("D8").Value = "'0.3810"
$ws.Range("E8").Value = "  -0.58%  "

$ws.Range("D9").Value = "'0.07748"
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").Value = "'0.9766"
$ws.Range("E10").Value = "  -1.14%  "

$ws.Range("D11").Value = "'22.61"
$ws.Range("E11").Value = "  +2.30%  "

$ws.Range("D12").Value = "1.946.48"
$ws.Range("E12").Value = "  +2.60%  "

$ws.Range("D13").Value = "'5.709"
$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").Value = "'6.943"
$ws.Range("E14").Value = "  -0.88%  "

$ws.Range("D15").Value = "'0.07016"
$ws.Range("E15").Value = "  -0.51%  "

$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").Value = "'84.62"
$ws.Range("E17").Value = "  +0.52%  "

$ws.Range("D18").Value = "'0.000009490"
$ws.Range("E18").Value = "  -0.66%  "

$ws.Range("D19").Value = "'16.67"
$ws.Range("E19").Value = "  -0.57%  "

$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("D21").Value = "29.023.48"
$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").Value = "'5.349"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").Value = "'11.03"
$ws.Range("E23").Value = "  +0.42%  "

$ws.Range("D24").Value = "2.154.71"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("D25").Value = "'2.061"
$ws.Range("E25").Value = "  -0.69%  "

$ws.Range("D26").Value = "'157.88"
$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("D27").Value = "'19.04"
$ws.Range("E27").Value = "  -1.08%  "

$ws.Range("D28").Value = "'5.588"
$ws.Range("E28").Value = "  -0.42%  "

$ws.Range("D29").Value = "'117.45"
$ws.Range("E29").Value = "  -0.49%  "

$ws.Range("D30").Value = "'1.829"
$ws.Range("E30").Value = "  -0.16%  "

$ws.Range("D31").Value = "'0.09323"
$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("D32").Value = "'0.8577"
$ws.Range("E32").Value = "  -0.65%  "

$ws.Range("D33").Value = "'5.084"
$ws.Range("E33").Value = "  -0.68%  "

$ws.Range("D34").Value = "'1.238"
$ws.Range("E34").Value = "  -1.49%  "

$ws.Range("D35").Value = "'3.022"
$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("E36").Value = "  -1.03%  "

$ws.Range("E37").Value = "  +0.38%  "

$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02039"
$ws.Range("E39").Value = "  -0.50%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'3.117"
$ws.Range("E40").Value = "  +13.76%  "

$ws.Range("D41").Value = "'7.420"
$ws.Range("E41").Value = "  -0.88%  "

$ws.Range("D42").Value = "'0.5486"
$ws.Range("E42").Value = "  -0.94%  "

$ws.Range("D43").Value = "'0.1753"
$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("D44").Value = "'9.341"
$ws.Range("E44").Value = "  +0.65%  "

$ws.Range("D45").Value = "'0.000002843"
$ws.Range("E45").Value = "  +12.72%  "

$ws.Range("D46").Value = "'2.175"
$ws.Range("E46").Value = "  +4.15%  "

$ws.Range("D47").Value = "'0.5176"
$ws.Range("E47").Value = "  -0.93%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.06917"
$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'11.18"
$ws.Range("E49").Value = "  -0.94%  "

$ws.Range("D50").Value = "'110.18"
$ws.Range("E50").Value = "  -1.79%  "

$ws.Range("D51").Value = "'1.755"
$ws.Range("E51").Value = "  -1.66%  "
